# Update column G ("K") on Sheet1 with newly calculated strikeout values,
# replacing the previously stored "Strike#" values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 2
    3  = 4
    4  = 3
    5  = 3
    6  = 3
    7  = 7
    8  = 5
    9  = 3
    10 = 3
    11 = 6
    12 = 3
    13 = 2
    14 = 3
    15 = 2
    16 = 2
    17 = 3
    18 = 5
    19 = 1
    20 = 3
    21 = 6
    22 = 1
    23 = 3
    24 = 2
    25 = 4
    26 = 3
    27 = 6
    28 = 0
    29 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

$wb.Save()
